$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1886711299453826
$ws.Range("C2").Value = 1.017921723714049
$ws.Range("D2").Value = 4.266585331249175
$ws.Range("E2").Value = 2.065571429713621
$ws.Range("F2").Value = 2.081870163654579
$ws.Range("G2").Value = 42

# Row 3
$ws.Range("B3").Value = 0.1733767525321955
$ws.Range("C3").Value = 1.062798336463839
$ws.Range("D3").Value = 4.575240169899562
$ws.Range("E3").Value = 2.138981105549921
$ws.Range("F3").Value = 2.158427712102139
$ws.Range("G3").Value = 41

# Row 4
$ws.Range("B4").Value = 0.1743475823563
$ws.Range("C4").Value = 1.003154296213989
$ws.Range("D4").Value = 4.399256932327289
$ws.Range("E4").Value = 2.097440567054831
$ws.Range("F4").Value = 2.116809367978816
$ws.Range("G4").Value = 40

# Row 5
$ws.Range("B5").Value = 0.1961419244217442
$ws.Range("C5").Value = 1.097619827889188
$ws.Range("D5").Value = 4.801704910896162
$ws.Range("E5").Value = 2.191279286375007
$ws.Range("F5").Value = 2.211013681542756
$ws.Range("G5").Value = 39

# Row 6
$ws.Range("B6").Value = 0.161080004310003
$ws.Range("C6").Value = 1.037884436753809
$ws.Range("D6").Value = 4.647468959415798
$ws.Range("E6").Value = 2.155798914420313
$ws.Range("F6").Value = 2.178629889817544
$ws.Range("G6").Value = 38

# Row 7
$ws.Range("B7").Value = 0.1938794397532232
$ws.Range("C7").Value = 1.078064803072593
$ws.Range("D7").Value = 4.91315197140037
$ws.Range("E7").Value = 2.216563098898917
$ws.Range("F7").Value = 2.238525191373713
$ws.Range("G7").Value = 37

# Row 8
$ws.Range("B8").Value = 0.1227245904408841
$ws.Range("C8").Value = 0.9933980158749599
$ws.Range("D8").Value = 4.520655172132273
$ws.Range("E8").Value = 2.126183240488052
$ws.Range("F8").Value = 2.152748266706019
$ws.Range("G8").Value = 36

# Row 9
$ws.Range("B9").Value = 0.1006825760723657
$ws.Range("C9").Value = 0.9849251405554545
$ws.Range("D9").Value = 4.66204020979998
$ws.Range("E9").Value = 2.159175817250643
$ws.Range("F9").Value = 2.188315313630956
$ws.Range("G9").Value = 35

# Row 10
$ws.Range("B10").Value = 0.1565014654937074
$ws.Range("C10").Value = 0.9857090473466144
$ws.Range("D10").Value = 4.70445062118944
$ws.Range("E10").Value = 2.168974555219457
$ws.Range("F10").Value = 2.19585400674243
$ws.Range("G10").Value = 34

# Row 11
$ws.Range("B11").Value = 0.1457776698747011
$ws.Range("C11").Value = 0.9842029885340268
$ws.Range("D11").Value = 4.813624961259621
$ws.Range("E11").Value = 2.19399748433302
$ws.Range("F11").Value = 2.223091431876469
$ws.Range("G11").Value = 33
